$d = $word.ActiveDocument
$ns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- Step 1: fix "NEXT: Serializ" + bookmark + "ation and Computer Strategy" -> one run ---
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $p.Range
$r.Find.Execute("NEXT: Serializ" + "ation and Computer Strategy", $true, $false, $false, $false, $false, $true, 1, $false, "NEXT: Serialization and Computer Strategy", 2) | Out-Null

# --- Step 2: build a buffer of empty paragraphs after the "NEXT: ..." paragraph so that
#             InsertXML never lands on the ambiguous "end of document" position. ---
$lastIdx = $d.Paragraphs.Count
$tail = $d.Paragraphs.Item($lastIdx)
for ($i = 0; $i -lt 7; $i++) {
    $tail.Range.InsertParagraphAfter() | Out-Null
    $tail = $d.Paragraphs.Item($d.Paragraphs.Count)
}

$idx = $lastIdx + 1

# Paragraph: blank spacer before the date line
$p1 = $d.Paragraphs.Item($idx)
$xml1 = "<w:p $ns><w:pPr><w:pStyle w:val=`"ListBullet`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"0`"/></w:numPr></w:pPr></w:p>"
$p1.Range.InsertXML($xml1) | Out-Null
$p1b = $d.Paragraphs.Item($idx)
$p1b.LeftIndent = 18
$p1b.FirstLineIndent = -18
$idx++

# Paragraph: bold date header "09/21/2016 (3hr)"
$p2 = $d.Paragraphs.Item($idx)
$xml2 = "<w:p $ns><w:pPr><w:pStyle w:val=`"ListBullet`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"0`"/></w:numPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>09/21/2016 (3hr)</w:t></w:r></w:p>"
$p2.Range.InsertXML($xml2) | Out-Null
$p2b = $d.Paragraphs.Item($idx)
$p2b.LeftIndent = 18
$p2b.FirstLineIndent = -18
$idx++

# Paragraph: "Implemented first half of serialization in BoardView class..."
$p3 = $d.Paragraphs.Item($idx)
$xml3 = "<w:p $ns><w:pPr><w:pStyle w:val=`"ListBullet`"/></w:pPr>" +
        "<w:r><w:t xml:space=`"preserve`">Implemented first half of serialization in </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>BoardView</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`"> class where the program asks the user their wish to serialize in every step along the game. Writes to file and exits the game if user wishes to.</w:t></w:r>" +
        "</w:p>"
$p3.Range.InsertXML($xml3) | Out-Null
$idx++

# Paragraph: "Found bug - If user entered more parameters than necessary..."
$p4 = $d.Paragraphs.Item($idx)
$xml4 = "<w:p $ns><w:pPr><w:pStyle w:val=`"ListBullet`"/></w:pPr>" +
        "<w:r><w:t xml:space=`"preserve`">Found bug " + [char]0x2013 + " If user entered more parameters than necessary, they were being carried over to the next </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>cin</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`"> request causing discrepancies in values assignment to variables.</w:t></w:r>" +
        "</w:p>"
$p4.Range.InsertXML($xml4) | Out-Null
$idx++

# Paragraph: "Fix: Added Cin.ignore() throughout the cins sections of code."
$p5 = $d.Paragraphs.Item($idx)
$xml5 = "<w:p $ns><w:pPr><w:pStyle w:val=`"ListBullet`"/>" +
        "<w:tabs><w:tab w:val=`"clear`" w:pos=`"360`"/><w:tab w:val=`"num`" w:pos=`"720`"/></w:tabs>" +
        "<w:ind w:left=`"720`"/></w:pPr>" +
        "<w:r><w:t xml:space=`"preserve`">Fix: Added </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/><w:proofErr w:type=`"gramStart`"/><w:r><w:t>Cin.ignore</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t>(</w:t></w:r><w:proofErr w:type=`"gramEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`">) throughout the </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>cins</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`"> sections of code.</w:t></w:r>" +
        "</w:p>"
$p5.Range.InsertXML($xml5) | Out-Null
$idx++

# Paragraph: "Next Step: Restoring the serialization file to the gameboard in the beginning of a game" + _GoBack bookmark
$p6 = $d.Paragraphs.Item($idx)
$xml6 = "<w:p $ns><w:pPr><w:pStyle w:val=`"ListBullet`"/></w:pPr>" +
        "<w:r><w:t xml:space=`"preserve`">Next Step: Restoring the serialization file to the </w:t></w:r>" +
        "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>gameboard</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
        "<w:r><w:t xml:space=`"preserve`"> in the beginning of a game</w:t></w:r>" +
        "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" +
        "</w:p>"
$p6.Range.InsertXML($xml6) | Out-Null
$idx++

# --- Step 3: remove the now-unneeded trailing blank paragraph buffer (keep doc length exact) ---
while ($d.Paragraphs.Count -ge $idx) {
    $last = $d.Paragraphs.Item($d.Paragraphs.Count)
    $t = $last.Range.Text
    if ($t.Length -le 1) {
        $last.Range.Delete() | Out-Null
    } else {
        break
    }
}

Write-Output "Done. Final paragraph count: $($d.Paragraphs.Count)"
